# Implemented writer to excel and fix known issue in 3.2
#
# This script:
#  - hides the (now-obsolete) "dataset" field rows (2-16)
#  - un-hides the "parameters" field rows (32-37), which are the newly
#    relevant rows after switching the sheet's active filter
#  - fixes row 37's "mandatory/optional" cell (E37) from "optional" to
#    "mandatory"
#  - renames the "uncertainty" -> "mean value" field (B53) to "mean"
#  - moves the sheet selection to E32 (first cell of the now-visible block)
#  - re-applies the AutoFilter on column A so it filters for "parameters"
#    instead of "dataset"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hide rows 2-16 ("dataset" field block)
for ($r = 2; $r -le 16; $r++) {
    $ws.Rows.Item($r).Hidden = $true
}

# Un-hide rows 32-37 ("parameters" field block)
for ($r = 32; $r -le 37; $r++) {
    $ws.Rows.Item($r).Hidden = $false
}

# Row 37 / column E: "optional" -> "mandatory"
$ws.Range("E37").Value = "mandatory"

# Row 53 / column B: "mean value" -> "mean"
$ws.Range("B53").Value = "mean"

# Move the active selection to E32
$ws.Range("E32").Select() | Out-Null

# Re-apply the autofilter on column A (dataset -> parameters)
$ws.Range("A1:G62").AutoFilter(1, @("parameters"), 7) | Out-Null

Write-Host "done"
